$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.158.37"
$ws.Range("E2").Value = "  -0.53%  "

$ws.Range("D3").Value = "1.863.57"
$ws.Range("E3").Value = "  -0.46%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4681"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.64%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2863"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.81%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06489"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.48%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07757"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.62%  "

$ws.Range("D12").Value = "1.863.16"
$ws.Range("E12").Value = "  -0.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "93.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6837"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.70%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.054"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "268.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.09%  "

$ws.Range("D17").Value = "30.153.58"
$ws.Range("E17").Value = "  -0.52%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007639"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.98%  "

$ws.Range("D21").Value = "2.114.98"
$ws.Range("E21").Value = "  -0.05%  "

$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.152"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.103"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.90%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.346"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.34%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.892"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.43%  "

$ws.Range("E29").Value = "  -0.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09925"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.29%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.453"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.225"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.008"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.00%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04688"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.50%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.118"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.74%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6896"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.707"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.22%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01836"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.76%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.757"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.94%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.356"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "71.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("D42").Style = "Normal"

$ws.Range("E43").Value = "  -3.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8342"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.85%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4061"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "933.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.128"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.70%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.958"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05578"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.87%  "
